# Commit message: "Added BAARD stage 2,  line, watermark to CIFAR10"
# The underlying XML diff shows that the second data row of the "apgd2"
# block (row 8, Adv_param = 0.2) was removed entirely, with every row
# below it (9-21) shifting up by one. This is exactly what Excel's
# "Delete Entire Row" does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("htru2")

# Delete row 8 completely - all rows below shift up, dimension
# shrinks from J21 to J20, and merged cell ranges (A8:A12, A14:A16,
# A18:A21) automatically shrink/shift to (A8:A11, A13:A15, A17:A20).
$ws.Rows.Item(8).Delete()

$wb.Save()
